$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("Testcases")
$wsData = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1. "Data" sheet: Runmode column (A) - all Chrome test rows should run,
#    flip the remaining "N" rows to "Y" (rows 4-13).
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 13; $r++) {
    $wsData.Cells.Item($r, 1).Value = "Y"
}

# ---------------------------------------------------------------------------
# 2. "Data" sheet: new "Total Jobs" numeric column (F) gets a number format.
#    Re-create the per-region styles (header / data / computed cell) instead
#    of letting Excel merge the format into the existing border/fill xf, so
#    the resulting style table mirrors the original edit (separate xfs for
#    header vs body vs the wrapped total cell).
# ---------------------------------------------------------------------------

# Header (F2) - same fill/border as the other header cells, plus numFmt "0"
$wsData.Cells.Item(2, 2).Copy() | Out-Null
$wsData.Cells.Item(2, 6).PasteSpecial(-4122) | Out-Null
$wsData.Cells.Item(2, 6).NumberFormat = "0"

# Body (F4:F13) - same border-only style as the other data columns, plus numFmt "0"
$wsData.Cells.Item(4, 4).Copy() | Out-Null
$bodyRng = $wsData.Range($wsData.Cells.Item(4, 6), $wsData.Cells.Item(13, 6))
$bodyRng.PasteSpecial(-4122) | Out-Null
$bodyRng.NumberFormat = "0"

# Computed total (F3) - plain numFmt "0" cell with wrapped text, holding the
# job count result.
$f3 = $wsData.Cells.Item(3, 6)
$f3.ClearFormats() | Out-Null
$f3.NumberFormat = "0"
$f3.WrapText = $true
$f3.Value = 4055

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Restore selections / active sheet to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsTest.Activate() | Out-Null
$wsTest.Range("A1:B3").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("C21").Select() | Out-Null
